$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: merge "DoS: " + "Cybersecurity Incident Report" into a single run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "DoS: Cybersecurity Incident Report", $true, $false, $false, $false, $false,
    $true, 1, $false, "DoS: Cybersecurity Incident Report", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Body paragraph: remove first-person language around "tcpdump"
#    "After using the network protocol analyzer tool, tcpdump, I have
#    confirmed ..." ->
#    "After using the network protocol analyzer tool, tcpdump, the
#    cybersecurity analyst has confirmed ..."
#
#    The original run layout is:
#       RunA("After using the network protocol analyzer tool, ")
#       proofErr(spellStart)
#       RunB("tcpdump")
#       proofErr(spellEnd)
#       RunC(", I have confirmed that UDP port 53 ... restore access to the site. ")
#
#    Doing the text edit in two Find/Replace passes that each straddle a
#    proofErr boundary folds both stray spell-check markers away (the
#    replacement text spans and subsumes the runs on both sides of each
#    marker). That leaves one merged run; re-toggling Bold off/on/off across
#    the "the cybersecurity analyst has" span then forces Word to re-split
#    that merged run back into three runs with identical (inherited)
#    formatting, matching the target shape.
# ---------------------------------------------------------------------------

# Pass 1: extend RunA's text through RunB, eating the spellStart marker.
$d.Content.Find.Execute(
    "tool, tcpdump", $true, $false, $false, $false, $false,
    $true, 1, $false, "tool, tcpdump, the cybersecurity analyst has", 2) | Out-Null

# Pass 2: join the new tail into RunC's text, eating the spellEnd marker.
$d.Content.Find.Execute(
    "has, I have confirmed", $true, $false, $false, $false, $false,
    $true, 1, $false, "has confirmed", 2) | Out-Null

# Re-split the merged run into three runs at the desired boundaries by
# toggling a character-formatting property across just the middle span and
# reverting it, without leaving any residual formatting difference.
$rng = $d.Content
$rng.Find.Execute("the cybersecurity analyst has") | Out-Null
$rng.Bold = 1
$rng.Bold = 0
